$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 281
$ws.Range("I2").Value = 97.75
$ws.Range("K2").Value = 97.75
$ws.Range("M2").Value = 15.25
$ws.Range("H17").Value = 1494581.6
$ws.Range("J17").Value = 1540547.8
$ws.Range("L17").Value = 4621643.4
$ws.Range("N17").Value = -4621979.4
$ws.Range("H18").Value = 169.5
$ws.Range("I18").Value = 169.5
$ws.Range("K18").Value = 169.5
$ws.Range("M18").Value = 114.5
$ws.Range("H28").Value = 310.45456
$ws.Range("I28").Value = 157.22223
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 157.22223
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 327.77777
$ws.Range("N28").Value = -1970
$ws.Range("H62").Value = 2779.5
$ws.Range("I62").Value = 2338.6924
$ws.Range("J62").Value = 3416.2222
$ws.Range("K62").Value = 2338.6924
$ws.Range("L62").Value = 3416.2222
$ws.Range("M62").Value = -1714.6924
$ws.Range("N62").Value = -4664.2222
$ws.Range("H65").Value = 2779.5
$ws.Range("I65").Value = 2338.6924
$ws.Range("J65").Value = 3416.2222
$ws.Range("K65").Value = 11693.462
$ws.Range("L65").Value = 17081.111
$ws.Range("M65").Value = -8573.462
$ws.Range("N65").Value = -23321.111
$ws.Range("H86").Value = 6199.625
$ws.Range("I86").Value = 1686.0714
$ws.Range("J86").Value = 12518.6
$ws.Range("K86").Value = 1686.0714
$ws.Range("L86").Value = 12518.6
$ws.Range("M86").Value = -563.0714
$ws.Range("N86").Value = -14764.6
$ws.Range("H89").Value = 6199.625
$ws.Range("I89").Value = 1686.0714
$ws.Range("J89").Value = 12518.6
$ws.Range("K89").Value = 8430.357
$ws.Range("L89").Value = 62593
$ws.Range("M89").Value = -2814.357
$ws.Range("N89").Value = -73825
$ws.Range("H96").Value = 50000908
$ws.Range("I96").Value = 62500930
$ws.Range("J96").Value = 829
$ws.Range("K96").Value = 187502790
$ws.Range("L96").Value = 2487
$ws.Range("M96").Value = -187501417
$ws.Range("N96").Value = -5233
$ws.Range("H121").Value = 1808.75
$ws.Range("J121").Value = 2030.5
$ws.Range("L121").Value = 6091.5
$ws.Range("N121").Value = -9585.5
$ws.Range("H132").Value = 1848.3726
$ws.Range("I132").Value = 1845.34
$ws.Range("K132").Value = 5536.02
$ws.Range("M132").Value = -3006.02
$ws.Range("H137").Value = 2761.261
$ws.Range("I137").Value = 2668.8948
$ws.Range("K137").Value = 8006.6844
$ws.Range("M137").Value = -5456.6844
$ws.Range("H141").Value = 1701.6666
$ws.Range("I141").Value = 1492.7273
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 4478.1819
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 701.8181000000004
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4273.431
$ws.Range("I32").Value = 4330.463
$ws.Range("J32").Value = 3503.5
$ws.Range("K32").Value = 4330.463
$ws.Range("L32").Value = 3503.5
$ws.Range("M32").Value = -4043.463
$ws.Range("N32").Value = -4077.5
$ws.Range("H63").Value = 1998.75
$ws.Range("I63").Value = 1998.75
$ws.Range("K63").Value = 1998.75
$ws.Range("M63").Value = -1312.75
$ws.Range("H66").Value = 1998.75
$ws.Range("I66").Value = 1998.75
$ws.Range("K66").Value = 9993.75
$ws.Range("M66").Value = -6561.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2731.8462
$ws.Range("I20").Value = 4443.7144
$ws.Range("J20").Value = 734.6667
$ws.Range("K20").Value = 4443.7144
$ws.Range("L20").Value = 734.6667
$ws.Range("M20").Value = -4196.7144
$ws.Range("N20").Value = -1228.6667
$ws.Range("H94").Value = 2315.15
$ws.Range("I94").Value = 2112.6875
$ws.Range("J94").Value = 3125
$ws.Range("K94").Value = 2112.6875
$ws.Range("L94").Value = 3125
$ws.Range("M94").Value = -1661.6875
$ws.Range("N94").Value = -4027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13281.923
$ws.Range("I31").Value = 21530.05
$ws.Range("K31").Value = 21530.05
$ws.Range("M31").Value = -21235.05
$ws.Range("H34").Value = 13281.923
$ws.Range("I34").Value = 21530.05
$ws.Range("K34").Value = 21530.05
$ws.Range("M34").Value = -21328.05
$ws.Range("H80").Value = 21900
$ws.Range("J80").Value = 21900
$ws.Range("L80").Value = 21900
$ws.Range("N80").Value = -24146
$ws.Range("H83").Value = 21900
$ws.Range("J83").Value = 21900
$ws.Range("L83").Value = 65700
$ws.Range("N83").Value = -76932
$ws.Range("H122").Value = 959.6818
$ws.Range("I122").Value = 1021.2
$ws.Range("J122").Value = 908.4167
$ws.Range("K122").Value = 3063.6
$ws.Range("L122").Value = 2725.2501
$ws.Range("M122").Value = -613.6000000000004
$ws.Range("N122").Value = -7625.2501
$ws.Range("H134").Value = 1067.6307
$ws.Range("I134").Value = 797.12195
$ws.Range("J134").Value = 1529.75
$ws.Range("K134").Value = 2391.36585
$ws.Range("L134").Value = 4589.25
$ws.Range("M134").Value = 143.6341499999999
$ws.Range("N134").Value = -9659.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4443.778
$ws.Range("I70").Value = 3731
$ws.Range("K70").Value = 11193
$ws.Range("M70").Value = -10878
$ws.Range("H73").Value = 4443.778
$ws.Range("I73").Value = 3731
$ws.Range("K73").Value = 11193
$ws.Range("M73").Value = -10101
$ws.Range("H131").Value = 786.27
$ws.Range("J131").Value = 790.8484999999999
$ws.Range("L131").Value = 2372.5455
$ws.Range("N131").Value = -12452.5455

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4118765
$ws.Range("J7").Value = 3127375.5
$ws.Range("L7").Value = 3127375.5
$ws.Range("N7").Value = -3127599.5
$ws.Range("H8").Value = 4118765
$ws.Range("J8").Value = 3127375.5
$ws.Range("L8").Value = 3127375.5
$ws.Range("N8").Value = -3127653.5
$ws.Range("H80").Value = 3152.5356
$ws.Range("I80").Value = 3040.5
$ws.Range("J80").Value = 3214.7778
$ws.Range("K80").Value = 3040.5
$ws.Range("L80").Value = 3214.7778
$ws.Range("M80").Value = -2042.5
$ws.Range("N80").Value = -5210.7778
$ws.Range("H83").Value = 3152.5356
$ws.Range("I83").Value = 3040.5
$ws.Range("J83").Value = 3214.7778
$ws.Range("K83").Value = 15202.5
$ws.Range("L83").Value = 16073.889
$ws.Range("M83").Value = -10210.5
$ws.Range("N83").Value = -26057.889
$ws.Range("H132").Value = 18579.395
$ws.Range("I132").Value = 3358.9614
$ws.Range("K132").Value = 10076.8842
$ws.Range("M132").Value = -7546.8842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2815.0625
$ws.Range("I7").Value = 3734.3333
$ws.Range("J7").Value = 2263.5
$ws.Range("K7").Value = 3734.3333
$ws.Range("L7").Value = 2263.5
$ws.Range("M7").Value = -3622.3333
$ws.Range("N7").Value = -2487.5
$ws.Range("H16").Value = 528.5714
$ws.Range("I16").Value = 600
$ws.Range("K16").Value = 600
$ws.Range("M16").Value = -430
$ws.Range("H45").Value = 16682
$ws.Range("J45").Value = 16682
$ws.Range("L45").Value = 16682
$ws.Range("N45").Value = -17496
$ws.Range("H61").Value = 3599.6785
$ws.Range("I61").Value = 1909
$ws.Range("K61").Value = 1909
$ws.Range("M61").Value = -1707
$ws.Range("H68").Value = 2624.111
$ws.Range("I68").Value = 2348
$ws.Range("J68").Value = 2969.25
$ws.Range("K68").Value = 2348
$ws.Range("L68").Value = 2969.25
$ws.Range("M68").Value = -1599
$ws.Range("N68").Value = -4467.25
$ws.Range("H71").Value = 2624.111
$ws.Range("I71").Value = 2348
$ws.Range("J71").Value = 2969.25
$ws.Range("K71").Value = 11740
$ws.Range("L71").Value = 14846.25
$ws.Range("M71").Value = -7996
$ws.Range("N71").Value = -22334.25
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H113").Value = 3599.6785
$ws.Range("I113").Value = 1909
$ws.Range("K113").Value = 1909
$ws.Range("M113").Value = 261
$ws.Range("H126").Value = 2815.0625
$ws.Range("I126").Value = 3734.3333
$ws.Range("J126").Value = 2263.5
$ws.Range("K126").Value = 11202.9999
$ws.Range("L126").Value = 6790.5
$ws.Range("M126").Value = -8732.999899999999
$ws.Range("N126").Value = -11730.5
$ws.Range("H132").Value = 1773.2572
$ws.Range("I132").Value = 983.9259
$ws.Range("J132").Value = 4437.25
$ws.Range("K132").Value = 2951.7777
$ws.Range("L132").Value = 13311.75
$ws.Range("M132").Value = -421.7776999999996
$ws.Range("N132").Value = -18371.75
$ws.Range("H136").Value = 25997.25
$ws.Range("I136").Value = 30409.176
$ws.Range("K136").Value = 91227.52799999999
$ws.Range("M136").Value = -88677.52799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5022.615
$ws.Range("I62").Value = 4206.25
$ws.Range("J62").Value = 5385.4443
$ws.Range("K62").Value = 4206.25
$ws.Range("L62").Value = 5385.4443
$ws.Range("M62").Value = -3582.25
$ws.Range("N62").Value = -6633.4443
$ws.Range("H65").Value = 5022.615
$ws.Range("I65").Value = 4206.25
$ws.Range("J65").Value = 5385.4443
$ws.Range("K65").Value = 21031.25
$ws.Range("L65").Value = 26927.2215
$ws.Range("M65").Value = -17911.25
$ws.Range("N65").Value = -33167.2215
$ws.Range("H70").Value = 15950
$ws.Range("J70").Value = 15950
$ws.Range("L70").Value = 15950
$ws.Range("N70").Value = -16580
$ws.Range("H73").Value = 15950
$ws.Range("J73").Value = 15950
$ws.Range("L73").Value = 15950
$ws.Range("N73").Value = -18134
$ws.Range("H132").Value = 1369.3334
$ws.Range("I132").Value = 1053.7
$ws.Range("K132").Value = 3161.1
$ws.Range("M132").Value = -631.1000000000004
$ws.Range("H136").Value = 27779506
$ws.Range("I136").Value = 41668160
$ws.Range("J136").Value = 2196.1667
$ws.Range("K136").Value = 125004480
$ws.Range("L136").Value = 6588.500100000001
$ws.Range("M136").Value = -125001930
$ws.Range("N136").Value = -11688.5001

